# fix units of nconc obs data
# The NConc (nitrogen concentration) columns (X, Y, Z, AA) were recorded
# as percentages but should be stored as fractions (i.e. divide by 100).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Observed")

$nconcCells = @{
    "X18"  = 6.32
    "Y18"  = 5.16
    "X19"  = 6.66
    "Y19"  = 5.77
    "X20"  = 6.53
    "Y20"  = 5.43
    "X21"  = 5.59
    "Y21"  = 2.83
    "AA21" = 3.24
    "X22"  = 6.17
    "Y22"  = 3.82
    "AA22" = 3.8
    "X23"  = 5.54
    "Y23"  = 2.88
    "AA23" = 3.23
    "Z27"  = 1.25
    "Z28"  = 1.44
    "Z29"  = 1.02
}

foreach ($addr in $nconcCells.Keys) {
    $ws.Range($addr).Value = $nconcCells[$addr] / 100
}

# Remove the autofilter that was left over on the observed data range.
$ws.AutoFilterMode = $false

# Stamp the sheet as unofficial output via header/footer watermark text.
$ps = $ws.PageSetup
$ps.CenterHeader = '&"Aptos"&12&KFF0000 UNOFFICIAL&1#' + [char]13
$ps.CenterFooter = [char]13 + '&1#&"Aptos"&12&KFF0000 UNOFFICIAL'

# Leave the selection on the last cell that was touched.
$ws.Range("Z29").Select() | Out-Null
